# Update Vtn-Itga5 LR-pairs sheet with new TPM-derived values
# and add the "Resolving-Mac" target/sending cluster rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Vtn", "Itga5", "ECs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 45.90594266666667, 137.717828, 0.3954672001633582, 0.3954672001633583, 127.0470834390187, 1143.423750951168, 0.01861263786402069, 0.0186126378640207),
    @("ECs", "Vtn", "Itga5", "FAPs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 40.44578266666667, 121.337348, 0.3484294080560655, 0.3484294080560656, 111.9358067106987, 1007.422260396288, 0.01639880726048523, 0.01639880726048523),
    @("ECs", "Vtn", "Itga5", "MuSCs", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 12.761795, 38.28538500000001, 0.1099393900775594, 0.1099393900775594, 35.31893127584, 317.87038148256, 0.005174290190588905, 0.005174290190588906),
    @("ECs", "Vtn", "Itga5", "Resolving-Mac", 3, 1, 2.767552, 8.302655999999999, 0.04706493447833917, 0.04706493447833917, 3, 1, 16.966758, 50.900274, 0.1461640017030168, 0.1461640017030168, 46.95638503641599, 422.6074653277439, 0.006879199163244341, 0.006879199163244342),
    @("FAPs", "Vtn", "Itga5", "ECs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 45.90594266666667, 137.717828, 0.3954672001633582, 0.3954672001633583, 966.4657679913956, 8698.19191192256, 0.14158906179244, 0.14158906179244),
    @("FAPs", "Vtn", "Itga5", "FAPs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 40.44578266666667, 121.337348, 0.3484294080560655, 0.3484294080560656, 851.5120730836622, 7663.60865775296, 0.1247481282067765, 0.1247481282067765),
    @("FAPs", "Vtn", "Itga5", "MuSCs", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 12.761795, 38.28538500000001, 0.1099393900775594, 0.1099393900775594, 268.6762821794667, 2418.0865396152, 0.03936158318233391, 0.03936158318233392),
    @("FAPs", "Vtn", "Itga5", "Resolving-Mac", 3, 1, 21.05317333333333, 63.15952, 0.3580298485789791, 0.3580298485789791, 3, 1, 16.966758, 50.900274, 0.1461640017030168, 0.1461640017030168, 357.20409707872, 3214.83687370848, 0.05233107539742875, 0.05233107539742876),
    @("MuSCs", "Vtn", "Itga5", "ECs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 45.90594266666667, 137.717828, 0.3954672001633582, 0.3954672001633583, 1605.671100504341, 14451.03990453907, 0.2352338512104116, 0.2352338512104117),
    @("MuSCs", "Vtn", "Itga5", "FAPs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 40.44578266666667, 121.337348, 0.3484294080560655, 0.3484294080560656, 1414.688830958314, 12732.19947862483, 0.2072545877335354, 0.2072545877335354),
    @("MuSCs", "Vtn", "Itga5", "MuSCs", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 12.761795, 38.28538500000001, 0.1099393900775594, 0.1099393900775594, 446.3745700824035, 4017.371130741631, 0.06539471823955376, 0.06539471823955377),
    @("MuSCs", "Vtn", "Itga5", "Resolving-Mac", 3, 1, 34.97741266666667, 104.932238, 0.5948251867999219, 0.5948251867999219, 3, 1, 16.966758, 50.900274, 0.1461640017030168, 0.1461640017030168, 593.4532961814681, 5341.079665633212, 0.08694202961642107, 0.08694202961642108),
    @("Resolving-Mac", "Vtn", "Itga5", "ECs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 45.90594266666667, 137.717828, 0.3954672001633582, 0.3954672001633583, 0.2160333661893334, 1.944300295704, 0.00003164929648587683, 0.00003164929648587684),
    @("Resolving-Mac", "Vtn", "Itga5", "FAPs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 40.44578266666667, 121.337348, 0.3484294080560655, 0.3484294080560656, 0.1903378532293333, 1.713040679064, 0.00002788485526842621, 0.00002788485526842621),
    @("Resolving-Mac", "Vtn", "Itga5", "MuSCs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 12.761795, 38.28538500000001, 0.1099393900775594, 0.1099393900775594, 0.06005700727, 0.5405130654300001, 0.000008798465082828214, 0.000008798465082828216),
    @("Resolving-Mac", "Vtn", "Itga5", "Resolving-Mac", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00008003014275976175, 0.00008003014275976175, 3, 1, 16.966758, 50.900274, 0.1461640017030168, 0.1461640017030168, 0.079845563148, 0.7186100683319999, 0.00001169752592263049, 0.0000116975259226305)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
